$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Update the "Currency" label/value row on the input sheet.
$ws1.Range("A6").Value = "currency"
$ws1.Range("B6").Value = "US Dollar"

# Restyle the row to the plain (non "Normal 2") look used elsewhere on the sheet.
$ws1.Range("A6").Interior.ThemeColor = 0
$ws1.Range("A6").Interior.TintAndShade = -0.34998626667073579
$ws1.Range("B6").Interior.Color = 5296274

# Make ProductLoanInput the active sheet/tab again (it was ProductLoanOutput),
# with the currency row selected.
$ws1.Activate()
$ws1.Range("A6:B6").Select()
